# Auto-generated Excel COM-interop script applying the Cactuar_Profits.xlsx diff
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H, I, J, K, L, M, N)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets to match the target commit.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 20543.908
$ws.Range("J40").Value = 22696.6
$ws.Range("L40").Value = 22696.6
$ws.Range("N40").Value = -23046.6
$ws.Range("H74").Value = 3483.4614
$ws.Range("I74").Value = 3830.8333
$ws.Range("J74").Value = 3185.7144
$ws.Range("K74").Value = 3830.8333
$ws.Range("L74").Value = 3185.7144
$ws.Range("M74").Value = -2894.8333
$ws.Range("N74").Value = -5057.7144
$ws.Range("H77").Value = 3483.4614
$ws.Range("I77").Value = 3830.8333
$ws.Range("J77").Value = 3185.7144
$ws.Range("K77").Value = 19154.1665
$ws.Range("L77").Value = 15928.572
$ws.Range("M77").Value = -14474.1665
$ws.Range("N77").Value = -25288.572
$ws.Range("H112").Value = 3066.7222
$ws.Range("I112").Value = 499
$ws.Range("K112").Value = 1497
$ws.Range("M112").Value = -389
$ws.Range("H131").Value = 3210.4546
$ws.Range("J131").Value = 12500
$ws.Range("L131").Value = 37500
$ws.Range("N131").Value = -47580
$ws.Range("H137").Value = 8806690
$ws.Range("I137").Value = 419607.4
$ws.Range("J137").Value = 22226024
$ws.Range("K137").Value = 1258822.2
$ws.Range("L137").Value = 66678072
$ws.Range("M137").Value = -1256272.2
$ws.Range("N137").Value = -66683172

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 822
$ws.Range("I22").Value = 822
$ws.Range("K22").Value = 822
$ws.Range("M22").Value = -523
$ws.Range("H32").Value = 11979.333
$ws.Range("I32").Value = 11382.933
$ws.Range("K32").Value = 11382.933
$ws.Range("M32").Value = -11095.933
$ws.Range("H45").Value = 105180.18
$ws.Range("I45").Value = 162711.72
$ws.Range("K45").Value = 162711.72
$ws.Range("M45").Value = -162334.72
$ws.Range("H63").Value = 1366.6666
$ws.Range("I63").Value = 1366.6666
$ws.Range("K63").Value = 1366.6666
$ws.Range("M63").Value = -680.6666
$ws.Range("H66").Value = 1366.6666
$ws.Range("I66").Value = 1366.6666
$ws.Range("K66").Value = 6833.333000000001
$ws.Range("M66").Value = -3401.333000000001
$ws.Range("H112").Value = 56666.332
$ws.Range("J112").Value = 56666.332
$ws.Range("L112").Value = 56666.332
$ws.Range("N112").Value = -59620.332
$ws.Range("H122").Value = 6184.857
$ws.Range("I122").Value = 4324
$ws.Range("J122").Value = 8666
$ws.Range("K122").Value = 12972
$ws.Range("L122").Value = 25998
$ws.Range("M122").Value = -10522
$ws.Range("N122").Value = -30898
$ws.Range("H132").Value = 31014.047
$ws.Range("I132").Value = 42867.855
$ws.Range("K132").Value = 128603.565
$ws.Range("M132").Value = -126073.565

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 37878.5
$ws.Range("I82").Value = 27838
$ws.Range("J82").Value = 68000
$ws.Range("K82").Value = 27838
$ws.Range("L82").Value = 68000
$ws.Range("M82").Value = -27455
$ws.Range("N82").Value = -68766
$ws.Range("H85").Value = 37878.5
$ws.Range("I85").Value = 27838
$ws.Range("J85").Value = 68000
$ws.Range("K85").Value = 27838
$ws.Range("L85").Value = 68000
$ws.Range("M85").Value = -26512
$ws.Range("N85").Value = -70652
$ws.Range("H108").Value = 60500
$ws.Range("I108").Value = 35000
$ws.Range("K108").Value = 35000
$ws.Range("M108").Value = -31160
$ws.Range("H125").Value = 21692.309
$ws.Range("J125").Value = 21692.309
$ws.Range("L125").Value = 21692.309
$ws.Range("N125").Value = -31532.309
$ws.Range("H134").Value = 2998.383
$ws.Range("I134").Value = 2967.96
$ws.Range("J134").Value = 3032.9546
$ws.Range("K134").Value = 8903.880000000001
$ws.Range("L134").Value = 9098.863799999999
$ws.Range("M134").Value = -6368.880000000001
$ws.Range("N134").Value = -14168.8638

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("J31").Value = 7611.8335
$ws.Range("L31").Value = 7611.8335
$ws.Range("N31").Value = -8201.833500000001
$ws.Range("J34").Value = 7611.8335
$ws.Range("L34").Value = 7611.8335
$ws.Range("N34").Value = -8015.8335
$ws.Range("H42").Value = 13000
$ws.Range("J42").Value = 12000.5
$ws.Range("L42").Value = 12000.5
$ws.Range("N42").Value = -13186.5
$ws.Range("H55").Value = 24877.555
$ws.Range("J55").Value = 31700
$ws.Range("L55").Value = 31700
$ws.Range("N55").Value = -32330
$ws.Range("H63").Value = 69950
$ws.Range("J63").Value = 69950
$ws.Range("L63").Value = 69950
$ws.Range("N63").Value = -71322
$ws.Range("H66").Value = 69950
$ws.Range("J66").Value = 69950
$ws.Range("L66").Value = 209850
$ws.Range("N66").Value = -216714
$ws.Range("H107").Value = 820.2258
$ws.Range("I107").Value = 496.33334
$ws.Range("K107").Value = 496.33334
$ws.Range("M107").Value = 1423.66666
$ws.Range("H134").Value = 1814.303
$ws.Range("I134").Value = 1245.7142
$ws.Range("K134").Value = 3737.1426
$ws.Range("M134").Value = -1202.1426

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H111").Value = 8231.75
$ws.Range("I111").Value = 8231.75
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 24695.25
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -21628.25
$ws.Range("N111").ClearContents()
$ws.Range("H113").Value = 678.3570999999999
$ws.Range("I113").Value = 649.75
$ws.Range("J113").Value = 689.8
$ws.Range("K113").Value = 1949.25
$ws.Range("L113").Value = 2069.4
$ws.Range("M113").Value = 220.75
$ws.Range("N113").Value = -6409.4
$ws.Range("H129").Value = 2189.9285
$ws.Range("I129").Value = 2475
$ws.Range("J129").Value = 1976.125
$ws.Range("K129").Value = 7425
$ws.Range("L129").Value = 5928.375
$ws.Range("M129").Value = -2425
$ws.Range("N129").Value = -15928.375
$ws.Range("H131").Value = 8773841
$ws.Range("I131").Value = 868
$ws.Range("J131").Value = 10754835
$ws.Range("K131").Value = 2604
$ws.Range("L131").Value = 32264505
$ws.Range("M131").Value = 2436
$ws.Range("N131").Value = -32274585
$ws.Range("H134").Value = 15591.333
$ws.Range("I134").Value = 8387
$ws.Range("K134").Value = 25161
$ws.Range("M134").Value = -20091

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 15875
$ws.Range("I57").Value = 2000
$ws.Range("J57").Value = 20500
$ws.Range("K57").Value = 2000
$ws.Range("L57").Value = 20500
$ws.Range("M57").Value = -1180
$ws.Range("N57").Value = -22140
$ws.Range("H102").Value = 11366122
$ws.Range("I102").Value = 13891269
$ws.Range("J102").Value = 2959.25
$ws.Range("K102").Value = 13891269
$ws.Range("L102").Value = 2959.25
$ws.Range("M102").Value = -13889647
$ws.Range("N102").Value = -6203.25
$ws.Range("H123").Value = 44452.8
$ws.Range("J123").Value = 44452.8
$ws.Range("L123").Value = 44452.8
$ws.Range("N123").Value = -49352.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 31750316
$ws.Range("I40").Value = 15876355
$ws.Range("J40").Value = 55561256
$ws.Range("K40").Value = 15876355
$ws.Range("L40").Value = 55561256
$ws.Range("M40").Value = -15876219
$ws.Range("N40").Value = -55561528
$ws.Range("H43").Value = 4896222
$ws.Range("J43").Value = 5505625
$ws.Range("L43").Value = 5505625
$ws.Range("N43").Value = -5506011
$ws.Range("H61").Value = 2223.7
$ws.Range("I61").Value = 1289.2858
$ws.Range("K61").Value = 1289.2858
$ws.Range("M61").Value = -1087.2858
$ws.Range("H68").Value = 4541.1665
$ws.Range("I68").Value = 3562
$ws.Range("J68").Value = 6499.5
$ws.Range("K68").Value = 3562
$ws.Range("L68").Value = 6499.5
$ws.Range("M68").Value = -2813
$ws.Range("N68").Value = -7997.5
$ws.Range("H71").Value = 4541.1665
$ws.Range("I71").Value = 3562
$ws.Range("J71").Value = 6499.5
$ws.Range("K71").Value = 17810
$ws.Range("L71").Value = 32497.5
$ws.Range("M71").Value = -14066
$ws.Range("N71").Value = -39985.5
$ws.Range("H93").Value = 1749.7894
$ws.Range("I93").Value = 1515.375
$ws.Range("J93").Value = 3000
$ws.Range("K93").Value = 1515.375
$ws.Range("L93").Value = 3000
$ws.Range("M93").Value = -267.375
$ws.Range("N93").Value = -5496
$ws.Range("H110").Value = 100000
$ws.Range("J110").Value = 100000
$ws.Range("L110").Value = 100000
$ws.Range("N110").Value = -108180
$ws.Range("H113").Value = 2223.7
$ws.Range("I113").Value = 1289.2858
$ws.Range("K113").Value = 1289.2858
$ws.Range("M113").Value = 880.7141999999999
$ws.Range("H132").Value = 2691.61
$ws.Range("I132").Value = 2640.6619
$ws.Range("K132").Value = 7921.9857
$ws.Range("M132").Value = -5391.9857
$ws.Range("H136").Value = 5081.8726
$ws.Range("I136").Value = 2949.9524
$ws.Range("J136").Value = 6398.647
$ws.Range("K136").Value = 8849.8572
$ws.Range("L136").Value = 19195.941
$ws.Range("M136").Value = -6299.8572
$ws.Range("N136").Value = -24295.941

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 179.09091
$ws.Range("I4").Value = 96.77778000000001
$ws.Range("K4").Value = 96.77778000000001
$ws.Range("M4").Value = 16.22221999999999
$ws.Range("H62").Value = 2290
$ws.Range("I62").Value = 2290
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 2290
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -1666
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 2290
$ws.Range("I65").Value = 2290
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 11450
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -8330
$ws.Range("N65").ClearContents()
$ws.Range("H132").Value = 4990.9546
$ws.Range("I132").Value = 1997.1818
$ws.Range("K132").Value = 5991.5454
$ws.Range("M132").Value = -3461.5454
